$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.91136961125256
$ws.Range("C2").Value = 6.981614341490489
$ws.Range("D2").Value = 8.659868089277644
$ws.Range("F2").Value = 35.90727507867637
$ws.Range("G2").Value = 3.679343829455548
$ws.Range("I2").Value = 26.85813808949183
$ws.Range("J2").Value = 10.39582709252655
$ws.Range("K2").Value = 10.94225883034519
$ws.Range("L2").Value = 11.75630123777942
$ws.Range("N2").Value = 20.38629481309523
$ws.Range("O2").Value = 27.39613999250987
# Row 3
$ws.Range("B3").Value = 13.70839579162905
$ws.Range("C3").Value = 6.936448089446534
$ws.Range("D3").Value = 8.64172854694894
$ws.Range("F3").Value = 35.98326052489513
$ws.Range("G3").Value = 3.681032852521445
$ws.Range("I3").Value = 26.94406105946747
$ws.Range("J3").Value = 10.41596661498348
$ws.Range("K3").Value = 10.80001463511535
$ws.Range("L3").Value = 11.7535638622998
$ws.Range("N3").Value = 20.44359653323766
$ws.Range("O3").Value = 27.47388480070775
# Row 4
$ws.Range("B4").Value = 13.58456563137349
$ws.Range("C4").Value = 6.908180257908715
$ws.Range("D4").Value = 8.631798343173685
$ws.Range("F4").Value = 36.03670866526616
$ws.Range("G4").Value = 3.682125962166448
$ws.Range("I4").Value = 27.0009468454777
$ws.Range("J4").Value = 10.42916497538205
$ws.Range("K4").Value = 10.71338430006119
$ws.Range("L4").Value = 11.75342990457175
$ws.Range("N4").Value = 20.48044839138423
$ws.Range("O4").Value = 27.52626241795823
# Row 5
$ws.Range("B5").Value = 13.53436474690903
$ws.Range("C5").Value = 6.896529907763189
$ws.Range("D5").Value = 8.628058406535445
$ws.Range("F5").Value = 36.06019588540175
$ws.Range("G5").Value = 3.68258554830043
$ws.Range("I5").Value = 27.02516658752718
$ws.Range("J5").Value = 10.43475323770334
$ws.Range("K5").Value = 10.67829899148495
$ws.Range("L5").Value = 11.75376574503383
$ws.Range("N5").Value = 20.49588655139375
$ws.Range("O5").Value = 27.54877331294812
# Row 6
$ws.Range("B6").Value = 13.52604645710725
$ws.Range("C6").Value = 6.894587573293776
$ws.Range("D6").Value = 8.62745600202533
$ws.Range("F6").Value = 36.06419895219889
$ws.Range("G6").Value = 3.682662717209504
$ws.Range("I6").Value = 27.02925097296477
$ws.Range("J6").Value = 10.43569385065227
$ws.Range("K6").Value = 10.67248736494763
$ws.Range("L6").Value = 11.75384513253643
$ws.Range("N6").Value = 20.49847549249932
$ws.Range("O6").Value = 27.55258166846299
# Row 7
$ws.Range("B7").Value = 13.58388746822095
$ws.Range("C7").Value = 6.908023662654483
$ws.Range("D7").Value = 8.631746659371728
$ws.Range("F7").Value = 36.03701851397483
$ws.Range("G7").Value = 3.682132103016817
$ws.Range("I7").Value = 27.00126927710204
$ws.Range("J7").Value = 10.42923949046623
$ws.Range("K7").Value = 10.71291019606398
$ws.Range("L7").Value = 11.75343285114709
$ws.Range("N7").Value = 20.48065489069838
$ws.Range("O7").Value = 27.52656128521133
# Row 8
$ws.Range("B8").Value = 13.84125523990152
$ws.Range("C8").Value = 6.966154383127029
$ws.Range("D8").Value = 8.653364841794183
$ws.Range("F8").Value = 35.93206423645744
$ws.Range("G8").Value = 3.679914598694308
$ws.Range("I8").Value = 26.88690748422499
$ws.Range("J8").Value = 10.40259866559775
$ws.Range("K8").Value = 10.8930896042223
$ws.Range("L8").Value = 11.75503727223867
$ws.Range("N8").Value = 20.40570690343543
$ws.Range("O8").Value = 27.42198256627653
# Row 9
$ws.Range("B9").Value = 14.34962299821024
$ws.Range("C9").Value = 7.075739609927647
$ws.Range("D9").Value = 8.705196931393292
$ws.Range("F9").Value = 35.78020239064401
$ws.Range("G9").Value = 3.676008805297509
$ws.Range("I9").Value = 26.69539852615581
$ws.Range("J9").Value = 10.35694200638184
$ws.Range("K9").Value = 11.25031593986738
$ws.Range("L9").Value = 11.77038778865523
$ws.Range("N9").Value = 20.27191608097911
$ws.Range("O9").Value = 27.25376509478324
# Row 10
$ws.Range("B10").Value = 14.7218571949703
$ws.Range("C10").Value = 7.153354863900767
$ws.Range("D10").Value = 8.748829417219758
$ws.Range("F10").Value = 35.70158424094242
$ws.Range("G10").Value = 3.673406391709079
$ws.Range("I10").Value = 26.57465685482084
$ws.Range("J10").Value = 10.32738480572525
$ws.Range("K10").Value = 11.51285035658728
$ws.Range("L10").Value = 11.78900716933601
$ws.Range("N10").Value = 20.18157661377062
$ws.Range("O10").Value = 27.15268486217098
# Row 11
$ws.Range("B11").Value = 14.89020307873332
$ws.Range("C11").Value = 7.187991187251966
$ws.Range("D11").Value = 8.769837083436206
$ws.Range("F11").Value = 35.67298236692418
$ws.Range("G11").Value = 3.672279921516101
$ws.Range("I11").Value = 26.5240609534031
$ws.Range("J11").Value = 10.31479819954526
$ws.Range("K11").Value = 11.63182500453523
$ws.Range("L11").Value = 11.79904794257799
$ws.Range("N11").Value = 20.14218962742365
$ws.Range("O11").Value = 27.1115955490005
# Row 12
$ws.Range("B12").Value = 14.95375186477323
$ws.Range("C12").Value = 7.201006820069772
$ws.Range("D12").Value = 8.777954476880282
$ws.Range("F12").Value = 35.66318151059627
$ws.Range("G12").Value = 3.671861564077904
$ws.Range("I12").Value = 26.50552423268855
$ws.Range("J12").Value = 10.31015506317778
$ws.Range("K12").Value = 11.67677368573119
$ws.Range("L12").Value = 11.80307371724642
$ws.Range("N12").Value = 20.12751928960585
$ws.Range("O12").Value = 27.09674015437534
# Row 13
$ws.Range("B13").Value = 14.94007530648635
$ws.Range("C13").Value = 7.198208196644925
$ws.Range("D13").Value = 8.776199107345342
$ws.Range("F13").Value = 35.66524648930944
$ws.Range("G13").Value = 3.671951300206158
$ws.Range("I13").Value = 26.50948874280485
$ws.Range("J13").Value = 10.31114957541604
$ws.Range("K13").Value = 11.66709844006621
$ws.Range("L13").Value = 11.80219679049895
$ws.Range("N13").Value = 20.13066794509928
$ws.Range("O13").Value = 27.09990819811882
# Row 14
$ws.Range("B14").Value = 14.89543561417428
$ws.Range("C14").Value = 7.18906402516373
$ws.Range("D14").Value = 8.770501680592725
$ws.Range("F14").Value = 35.67215540098287
$ws.Range("G14").Value = 3.672245338630899
$ws.Range("I14").Value = 26.52252344234003
$ws.Range("J14").Value = 10.31441374033339
$ws.Range("K14").Value = 11.63552527321696
$ws.Range("L14").Value = 11.79937467894063
$ws.Range("N14").Value = 20.14097779270107
$ws.Range("O14").Value = 27.11035926911369
# Row 15
$ws.Range("B15").Value = 14.86806474963009
$ws.Range("C15").Value = 7.183449768259237
$ws.Range("D15").Value = 8.767032837124274
$ws.Range("F15").Value = 35.67652145158137
$ws.Range("G15").Value = 3.672426514001634
$ws.Range("I15").Value = 26.5305886910586
$ws.Range("J15").Value = 10.3164291599838
$ws.Range("K15").Value = 11.61617105167668
$ws.Range("L15").Value = 11.79767509601156
$ws.Range("N15").Value = 20.14732470209931
$ws.Range("O15").Value = 27.11685258342589
# Row 16
$ws.Range("B16").Value = 14.71083052390839
$ws.Range("C16").Value = 7.15107745992937
$ws.Range("D16").Value = 8.747479465757896
$ws.Range("F16").Value = 35.70359757970609
$ws.Range("G16").Value = 3.673481160499929
$ws.Range("I16").Value = 26.57805056448481
$ws.Range("J16").Value = 10.32822462321275
$ws.Range("K16").Value = 11.50506257930183
$ws.Range("L16").Value = 11.78838239337956
$ws.Range("N16").Value = 20.1841849408649
$ws.Range("O16").Value = 27.15546865730827
# Row 17
$ws.Range("B17").Value = 14.61407781347673
$ws.Range("C17").Value = 7.131043409524771
$ws.Range("D17").Value = 8.735777874381014
$ws.Range("F17").Value = 35.72204248007197
$ws.Range("G17").Value = 3.674142820771344
$ws.Range("I17").Value = 26.60827605999261
$ws.Range("J17").Value = 10.33568051067306
$ws.Range("K17").Value = 11.43675671560808
$ws.Range("L17").Value = 11.78308231809443
$ws.Range("N17").Value = 20.20723442131873
$ws.Range("O17").Value = 27.18041199001879
# Row 18
$ws.Range("B18").Value = 14.55833884457043
$ws.Range("C18").Value = 7.119457401430567
$ws.Range("D18").Value = 8.729156798406629
$ws.Range("F18").Value = 35.73332562304856
$ws.Range("G18").Value = 3.674528793516766
$ws.Range("I18").Value = 26.62606852294994
$ws.Range("J18").Value = 10.34004982744924
$ws.Range("K18").Value = 11.39742857852223
$ws.Range("L18").Value = 11.78018184651121
$ws.Range("N18").Value = 20.22065278538367
$ws.Range("O18").Value = 27.19521921243133
# Row 19
$ws.Range("B19").Value = 14.53945305316003
$ws.Range("C19").Value = 7.11552388464185
$ws.Range("D19").Value = 8.726933931903162
$ws.Range("F19").Value = 35.73726167189762
$ws.Range("G19").Value = 3.674660406462253
$ws.Range("I19").Value = 26.63216273786032
$ws.Range("J19").Value = 10.34154310782234
$ws.Range("K19").Value = 11.38410701913352
$ws.Range("L19").Value = 11.77922528262151
$ws.Range("N19").Value = 20.22522368279247
$ws.Range("O19").Value = 27.20031174879808
# Row 20
$ws.Range("B20").Value = 14.62438698407637
$ws.Range("C20").Value = 7.13318261186846
$ws.Range("D20").Value = 8.737012240834076
$ws.Range("F20").Value = 35.72000921724987
$ws.Range("G20").Value = 3.674071826971983
$ws.Range("I20").Value = 26.60501632237349
$ws.Range("J20").Value = 10.33487845008918
$ws.Range("K20").Value = 11.4440324613436
$ws.Range("L20").Value = 11.78363121992254
$ws.Range("N20").Value = 20.2047641187639
$ws.Range("O20").Value = 27.17770906679367
# Row 21
$ws.Range("B21").Value = 14.90855326121537
$ws.Range("C21").Value = 7.191752645152001
$ws.Range("D21").Value = 8.77217078575692
$ws.Range("F21").Value = 35.67009813109535
$ws.Range("G21").Value = 3.672158749845401
$ws.Range("I21").Value = 26.51867792968635
$ws.Range("J21").Value = 10.31345163721607
$ws.Range("K21").Value = 11.64480221891058
$ws.Range("L21").Value = 11.80019755354264
$ws.Range("N21").Value = 20.1379429090994
$ws.Range("O21").Value = 27.10727041832631
# Row 22
$ws.Range("B22").Value = 15.09308016358895
$ws.Range("C22").Value = 7.229444632749734
$ws.Range("D22").Value = 8.796092445464103
$ws.Range("F22").Value = 35.64348212551581
$ws.Range("G22").Value = 3.670956293551681
$ws.Range("I22").Value = 26.46588132630225
$ws.Range("J22").Value = 10.30016553654484
$ws.Range("K22").Value = 11.77539077240409
$ws.Range("L22").Value = 11.81232646744865
$ws.Range("N22").Value = 20.09569695410918
$ws.Range("O22").Value = 27.06533978616517
# Row 23
$ws.Range("B23").Value = 14.99472252485056
$ws.Range("C23").Value = 7.209382677184248
$ws.Range("D23").Value = 8.783240180432793
$ws.Range("F23").Value = 35.65713827278666
$ws.Range("G23").Value = 3.671593701578345
$ws.Range("I23").Value = 26.4937276120097
$ws.Range("J23").Value = 10.3071910533292
$ws.Range("K23").Value = 11.70576314543703
$ws.Range("L23").Value = 11.805734719336
$ws.Range("N23").Value = 20.11811432886556
$ws.Range("O23").Value = 27.08734312875561
# Row 24
$ws.Range("B24").Value = 14.61972656035729
$ws.Range("C24").Value = 7.132215690530954
$ws.Range("D24").Value = 8.73645385224923
$ws.Range("F24").Value = 35.7209263405945
$ws.Range("G24").Value = 3.674103905898524
$ws.Range("I24").Value = 26.60648875574769
$ws.Range("J24").Value = 10.33524080365006
$ws.Range("K24").Value = 11.44074327492891
$ws.Range("L24").Value = 11.783382604556
$ws.Range("N24").Value = 20.20588042254375
$ws.Range("O24").Value = 27.17892960385803
# Row 25
$ws.Range("B25").Value = 14.21206871695207
$ws.Range("C25").Value = 7.046587595722931
$ws.Range("D25").Value = 8.690184585333824
$ws.Range("F25").Value = 35.81550124587973
$ws.Range("G25").Value = 3.6770183119588
$ws.Range("I25").Value = 26.74370104175076
$ws.Range("J25").Value = 10.36859120368289
$ws.Range("K25").Value = 11.15349832401352
$ws.Range("L25").Value = 11.77038778865523
$ws.Range("N25").Value = 20.30670707466863
$ws.Range("O25").Value = 27.29532202881195
